$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.978.95'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '1.744.66'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '249.78'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +6.79%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5138'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.32%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2761'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '1.742.12'
$ws.Range("E10").Value = '  -0.42%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07240'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '15.25'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.58%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.6500'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  +0.80%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '77.78'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("E16").Value = '  +0.05%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '26.008.72'
$ws.Range("E18").Value = '  +0.02%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.86'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.74%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000006820'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '1.962.75'
$ws.Range("E21").Value = '  -0.80%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.303'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.27%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.682'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.66%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.368'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +3.06%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '135.86'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.11%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.508'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("E27").Value = '  -0.51%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.789'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.64%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '106.22'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.92%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '3.945'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +4.15%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08253'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.26%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.670'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04679'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.69%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.654'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("E35").Value = '  -0.18%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6250'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("E37").Value = '  +0.18%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01606'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.25%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.933'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '100.67'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.84%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.3886'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.57%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.7568'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("E44").Value = '  -0.60%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '6.364'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("E48").Value = '  -2.19%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '30.72'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.61%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.580'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.89%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.3440'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
